$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3062.4375
$ws.Cells.Item(43, 9).Value = 3285.7144
$ws.Cells.Item(43, 10).Value = 1499.5
$ws.Cells.Item(43, 11).Value = 3285.7144
$ws.Cells.Item(43, 12).Value = 1499.5
$ws.Cells.Item(43, 13).Value = -3216.7144
$ws.Cells.Item(43, 14).Value = -1637.5

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 631.4

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1115.1904
$ws.Cells.Item(70, 9).Value = 994.9375
$ws.Cells.Item(70, 11).Value = 2984.8125
$ws.Cells.Item(70, 13).Value = -2714.8125

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1115.1904
$ws.Cells.Item(73, 9).Value = 994.9375
$ws.Cells.Item(73, 11).Value = 2984.8125
$ws.Cells.Item(73, 13).Value = -2048.8125

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3006.5483
$ws.Cells.Item(76, 9).Value = 3000.1
$ws.Cells.Item(76, 11).Value = 3000.1
$ws.Cells.Item(76, 13).Value = -2685.1

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 3006.5483
$ws.Cells.Item(79, 9).Value = 3000.1
$ws.Cells.Item(79, 11).Value = 3000.1
$ws.Cells.Item(79, 13).Value = -1908.1

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 22675.334
$ws.Cells.Item(88, 9).Value = 8210
$ws.Cells.Item(88, 11).Value = 8210
$ws.Cells.Item(88, 13).Value = -7804

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 22675.334
$ws.Cells.Item(91, 9).Value = 8210
$ws.Cells.Item(91, 11).Value = 8210
$ws.Cells.Item(91, 13).Value = -6806

# ALC row 114
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(114, 8).Value = 34992
$ws.Cells.Item(114, 10).Value = 34992
$ws.Cells.Item(114, 12).Value = 34992
$ws.Cells.Item(114, 14).Value = -43670

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value = 47745.6
$ws.Cells.Item(126, 10).Value = 47745.6
$ws.Cells.Item(126, 12).Value = 47745.6
$ws.Cells.Item(126, 14).Value = -57625.6

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(128, 8).Value = 52251.5
$ws.Cells.Item(128, 10).Value = 52251.5
$ws.Cells.Item(128, 12).Value = 52251.5
$ws.Cells.Item(128, 14).Value = -62211.5

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 301446.7
$ws.Cells.Item(129, 10).Value = 1331.7142
$ws.Cells.Item(129, 12).Value = 3995.1426
$ws.Cells.Item(129, 14).Value = -13995.1426

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17458.322
$ws.Cells.Item(32, 9).Value = 16437.135
$ws.Cells.Item(32, 11).Value = 16437.135
$ws.Cells.Item(32, 13).Value = -16150.135

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(109, 8).Value = 27821
$ws.Cells.Item(109, 10).Value = 27821
$ws.Cells.Item(109, 12).Value = 27821
$ws.Cells.Item(109, 14).Value = -30595

# ARM row 111
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(111, 8).Value = 41409.332
$ws.Cells.Item(111, 10).Value = 41409.332
$ws.Cells.Item(111, 12).Value = 41409.332
$ws.Cells.Item(111, 14).Value = -49589.332

# ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value = 51690
$ws.Cells.Item(119, 10).Value = 51690
$ws.Cells.Item(119, 12).Value = 51690
$ws.Cells.Item(119, 14).Value = -61366

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(121, 8).Value = 32498.5
$ws.Cells.Item(121, 10).Value = 32498.5
$ws.Cells.Item(121, 12).Value = 32498.5
$ws.Cells.Item(121, 14).Value = -35992.5

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(137, 8).Value = 59428
$ws.Cells.Item(137, 10).Value = 59428
$ws.Cells.Item(137, 12).Value = 59428
$ws.Cells.Item(137, 14).Value = -69628

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2389.2
$ws.Cells.Item(107, 9).Value = 2116.5
$ws.Cells.Item(107, 11).Value = 2116.5
$ws.Cells.Item(107, 13).Value = -196.5

# BSM row 110
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(110, 8).Value = 47030
$ws.Cells.Item(110, 10).Value = 47030
$ws.Cells.Item(110, 12).Value = 47030
$ws.Cells.Item(110, 14).Value = -55210

# BSM row 112
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 44600.25
$ws.Cells.Item(112, 10).Value = 44600.25
$ws.Cells.Item(112, 12).Value = 44600.25
$ws.Cells.Item(112, 14).Value = -47554.25

# BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(130, 8).Value = 36927
$ws.Cells.Item(130, 10).Value = 36927
$ws.Cells.Item(130, 12).Value = 36927
$ws.Cells.Item(130, 14).Value = -46967

# BSM row 139
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(139, 8).Value = 63123
$ws.Cells.Item(139, 10).Value = 63123
$ws.Cells.Item(139, 12).Value = 63123
$ws.Cells.Item(139, 14).Value = -73403

# CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(110, 8).Value = 40985.668
$ws.Cells.Item(110, 10).Value = 40985.668
$ws.Cells.Item(110, 12).Value = 40985.668
$ws.Cells.Item(110, 14).Value = -49165.668

# CRP row 112
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(112, 8).Value = 40492
$ws.Cells.Item(112, 10).Value = 40492
$ws.Cells.Item(112, 12).Value = 40492
$ws.Cells.Item(112, 14).Value = -43446

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 27278.215
$ws.Cells.Item(132, 9).Value = 1683.1471
$ws.Cells.Item(132, 11).Value = 5049.4413
$ws.Cells.Item(132, 13).Value = -2519.4413

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 46997.977
$ws.Cells.Item(131, 10).Value = 49035.906
$ws.Cells.Item(131, 12).Value = 147107.718
$ws.Cells.Item(131, 14).Value = -157187.718

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 45856.523
$ws.Cells.Item(140, 9).Value = 64481.25
$ws.Cells.Item(140, 10).Value = 3285.7144
$ws.Cells.Item(140, 11).Value = 193443.75
$ws.Cells.Item(140, 12).Value = 9857.143199999999
$ws.Cells.Item(140, 13).Value = -188263.75
$ws.Cells.Item(140, 14).Value = -20217.1432

# GSM row 7
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 100000
$ws.Cells.Item(7, 9).Value = 100000
$ws.Cells.Item(7, 11).Value = 100000
$ws.Cells.Item(7, 13).Value = -99888

# GSM row 8
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(8, 8).Value = 100000
$ws.Cells.Item(8, 9).Value = 100000
$ws.Cells.Item(8, 11).Value = 100000
$ws.Cells.Item(8, 13).Value = -99861

# GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(114, 8).Value = 44242
$ws.Cells.Item(114, 10).Value = 44242
$ws.Cells.Item(114, 12).Value = 44242
$ws.Cells.Item(114, 14).Value = -52920

# GSM row 116
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(116, 8).Value = 42346
$ws.Cells.Item(116, 10).Value = 42346
$ws.Cells.Item(116, 12).Value = 42346
$ws.Cells.Item(116, 14).Value = -51524

# GSM row 119
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(119, 8).Value = 47753
$ws.Cells.Item(119, 10).Value = 47753
$ws.Cells.Item(119, 12).Value = 47753
$ws.Cells.Item(119, 14).Value = -57429

# GSM row 128
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(128, 8).Value = 38758.668
$ws.Cells.Item(128, 10).Value = 38758.668
$ws.Cells.Item(128, 12).Value = 38758.668
$ws.Cells.Item(128, 14).Value = -48718.668

# LTW row 3
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()

# LTW row 15
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()

# LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(108, 8).Value = 48618
$ws.Cells.Item(108, 10).Value = 48618
$ws.Cells.Item(108, 12).Value = 48618
$ws.Cells.Item(108, 14).Value = -56298

# LTW row 112
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(112, 8).Value = 28938
$ws.Cells.Item(112, 10).Value = 28938
$ws.Cells.Item(112, 12).Value = 28938
$ws.Cells.Item(112, 14).Value = -31892

# LTW row 118
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(118, 8).Value = 43405
$ws.Cells.Item(118, 10).Value = 43405
$ws.Cells.Item(118, 12).Value = 43405
$ws.Cells.Item(118, 14).Value = -46719

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 40858.883
$ws.Cells.Item(122, 9).Value = 57462.832
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 172388.496
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -169938.496
$ws.Cells.Item(122, 14).Value = -15400

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 35141.715
$ws.Cells.Item(133, 10).Value = 35141.715
$ws.Cells.Item(133, 12).Value = 35141.715
$ws.Cells.Item(133, 14).Value = -40201.715

# LTW row 137
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(137, 8).Value = 59495
$ws.Cells.Item(137, 10).Value = 59495
$ws.Cells.Item(137, 12).Value = 59495
$ws.Cells.Item(137, 14).Value = -69695

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 44883
$ws.Cells.Item(139, 10).Value = 59499.5
$ws.Cells.Item(139, 12).Value = 59499.5
$ws.Cells.Item(139, 14).Value = -69779.5

# WVR row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 38686.5
$ws.Cells.Item(109, 10).Value = 38686.5
$ws.Cells.Item(109, 12).Value = 38686.5
$ws.Cells.Item(109, 14).Value = -41460.5

# WVR row 117
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(117, 8).Value = 50301
$ws.Cells.Item(117, 10).Value = 50301
$ws.Cells.Item(117, 12).Value = 50301
$ws.Cells.Item(117, 14).Value = -59479

# WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 39715
$ws.Cells.Item(125, 10).Value = 39715
$ws.Cells.Item(125, 12).Value = 39715
$ws.Cells.Item(125, 14).Value = -49555

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 60000
$ws.Cells.Item(139, 10).Value = 60000
$ws.Cells.Item(139, 12).Value = 60000
$ws.Cells.Item(139, 14).Value = -70280
